# Auto-generated: apply 2022-12-18 daily update to violent-crime-full-year workbook
# Each worksheet's column I (year 2022 running total) is bumped for the rows
# whose crime counts increased with the newly-added day of data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 7065
$ws.Range("I3").Value = 7284
$ws.Range("I4").Value = 1678
$ws.Range("I6").Value = 8647
$ws.Range("I7").Value = 25363

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 204
$ws.Range("I6").Value = 184
$ws.Range("I7").Value = 795
$ws.Range("I8").Value = 1506
$ws.Range("I10").Value = 186
$ws.Range("I11").Value = 388
$ws.Range("I18").Value = 200
$ws.Range("I19").Value = 707
$ws.Range("I20").Value = 626
$ws.Range("I21").Value = 113
$ws.Range("I29").Value = 1505
$ws.Range("I33").Value = 1118
$ws.Range("I34").Value = 115
$ws.Range("I36").Value = 348
$ws.Range("I37").Value = 783
$ws.Range("I40").Value = 45
$ws.Range("I42").Value = 955
$ws.Range("I50").Value = 134
$ws.Range("I52").Value = 575
$ws.Range("I53").Value = 290
$ws.Range("I54").Value = 496
$ws.Range("I55").Value = 299
$ws.Range("I56").Value = 28
$ws.Range("I63").Value = 77
$ws.Range("I65").Value = 590
$ws.Range("I66").Value = 76
$ws.Range("I67").Value = 960
$ws.Range("I73").Value = 229
$ws.Range("I76").Value = 363
$ws.Range("I78").Value = 336
$ws.Range("I79").Value = 728
$ws.Range("I83").Value = 548
$ws.Range("I85").Value = 1128
$ws.Range("I86").Value = 166
$ws.Range("I87").Value = 67
$ws.Range("I89").Value = 298
$ws.Range("I90").Value = 329
$ws.Range("I91").Value = 268
$ws.Range("I92").Value = 77
$ws.Range("I95").Value = 390
$ws.Range("I96").Value = 296
$ws.Range("I99").Value = 441
$ws.Range("I101").Value = 25363

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 323
$ws.Range("I3").Value = 423
$ws.Range("I7").Value = 1128

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 185
$ws.Range("I6").Value = 189
$ws.Range("I7").Value = 575

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 388

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 444
$ws.Range("I3").Value = 434
$ws.Range("I7").Value = 1506

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I6").Value = 144
$ws.Range("I7").Value = 290

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 261
$ws.Range("I4").Value = 42
$ws.Range("I6").Value = 216
$ws.Range("I7").Value = 795

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I4").Value = 42
$ws.Range("I6").Value = 105
$ws.Range("I7").Value = 298

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 83
$ws.Range("I7").Value = 296

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 232
$ws.Range("I3").Value = 254
$ws.Range("I6").Value = 238
$ws.Range("I7").Value = 783

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 125
$ws.Range("I6").Value = 115
$ws.Range("I7").Value = 441

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 231
$ws.Range("I7").Value = 960

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 177
$ws.Range("I6").Value = 180
$ws.Range("I7").Value = 590

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I3").Value = 200
$ws.Range("I7").Value = 548

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 139
$ws.Range("I7").Value = 390

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I6").Value = 359
$ws.Range("I7").Value = 1118

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 240
$ws.Range("I7").Value = 496

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I3").Value = 513
$ws.Range("I6").Value = 417
$ws.Range("I7").Value = 1505

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 228
$ws.Range("I3").Value = 205
$ws.Range("I7").Value = 707

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I6").Value = 167
$ws.Range("I7").Value = 363

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I3").Value = 46
$ws.Range("I7").Value = 184

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I6").Value = 385
$ws.Range("I7").Value = 955

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I6").Value = 88
$ws.Range("I7").Value = 186

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I2").Value = 76
$ws.Range("I7").Value = 336

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 299

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I2").Value = 84
$ws.Range("I7").Value = 268

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I3").Value = 17
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 213
$ws.Range("I3").Value = 239
$ws.Range("I6").Value = 206
$ws.Range("I7").Value = 728

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I4").Value = 40
$ws.Range("I7").Value = 626

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I3").Value = 45
$ws.Range("I7").Value = 200

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 99
$ws.Range("I6").Value = 108
$ws.Range("I7").Value = 348

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 115

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I3").Value = 68
$ws.Range("I7").Value = 229

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I2").Value = 77
$ws.Range("I7").Value = 204

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I2").Value = 27
$ws.Range("I7").Value = 77

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 80
$ws.Range("I7").Value = 166

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I3").Value = 86
$ws.Range("I7").Value = 329

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("I3").Value = 19
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 67
